# Daily attendance processing - 2026-01-06 01:38:43
# Swap the order of the two comma-separated "Recorded By" names in column G
# for every row where the list currently contains both "System" and the
# user's email address.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in column G (Recorded By) whose value needs its two comma-separated
# entries reversed in order, as identified from the diff.
$rows = @(8,9,10,12,14,15,17,18,34,35,36,38,40,41,43,44,60,61,62,64,66,67,69,70,86,87,88,90,92,93,95,96,112,113,114,116,118,119,121,122,138,139,140,142,144,145,147,148,164,167,170,174,191,194,197,201,218,221,224,228,245,248,251,255,272,275,278,282,299,302,305,309)

foreach ($row in $rows) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = 7
    $current = $cell.Value2

    if ($current -ne $null -and $current.ToString().Contains(",")) {
        $parts = $current.ToString().Split(",")
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }
        $reversed = $trimmed[($trimmed.Count - 1)..0]
        $cell.Value2 = [string]::Join(", ", $reversed)
    }
}
